# This script reproduces the diff:
#  - The original Sheet1 (square / location / loop_no table, A1:C13) is
#    moved to a new Sheet2 (same data, sheetView no longer tabSelected,
#    selection becomes A1:C13 with no explicit activeCell).
#  - Sheet1 is rewritten with the transposed version of that table:
#    row 1 = "square" header + the 12 square numbers,
#    row 2 = "location" header + the 12 location strings,
#    row 3 = "loop_no" header + the 12 loop numbers.
#    Sheet1 stays the tabSelected sheet, with selection on B19.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Read the existing Sheet1 data (3 columns x 13 rows: header + 12 rows) ---
$nRows = 13
$nCols = 3
$data = @()
for ($r = 1; $r -le $nRows; $r++) {
    $row = @()
    for ($c = 1; $c -le $nCols; $c++) {
        $row += $ws1.Cells.Item($r, $c).Value()
    }
    $data += ,$row
}

# --- Create Sheet2 (placed after Sheet1) and copy the original data into it ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

for ($r = 1; $r -le $nRows; $r++) {
    for ($c = 1; $c -le $nCols; $c++) {
        $ws2.Cells.Item($r, $c).Value = $data[$r - 1][$c - 1]
    }
}

# --- Clear Sheet1 and write the transposed table ---
$ws1.Cells.Clear()

for ($r = 1; $r -le $nRows; $r++) {
    for ($c = 1; $c -le $nCols; $c++) {
        # transposed: original row r, col c -> new row c, col r
        $ws1.Cells.Item($c, $r).Value = $data[$r - 1][$c - 1]
    }
}

# --- View / selection state ---
# Sheet2: no longer the active/selected tab; selection spans the whole table.
$ws2.Range("A1:C13").Select()

# Make Sheet1 active again and set its selection as in the target.
$ws1.Activate()
$ws1.Range("B19").Select()
